$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B -- shifts old B..K to C..L
$ws.Range("B1").EntireColumn.Insert()

# New header cell for the inserted column - copy the header style from a
# neighboring header cell, then set its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B1").Value = "segments"

# For each data row, move the text label currently in column A into column B,
# and replace column A with a zero-based numeric index.
for ($r = 2; $r -le 20; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $label
    $ws.Cells.Item($r, 2).Style = "Normal"
    $ws.Cells.Item($r, 1).Value = $r - 2
}
